$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Feb 15 12:43:41 EST 2023"
$ws.Range("B3").Value = "Wed Feb 15 12:43:51 EST 2023"
$ws.Range("B4").Value = "Wed Feb 15 12:44:00 EST 2023"
$ws.Range("B5").Value = "Wed Feb 15 12:44:10 EST 2023"
$ws.Range("B6").Value = "Wed Feb 15 12:44:19 EST 2023"
$ws.Range("B7").Value = "Wed Feb 15 12:44:29 EST 2023"
